# Fix Critical, High, and Medium QA issues: Config errors, CPO inflation, Prod zones,
# code cleanup, imports, and logging.
#
# This script:
#  1. Fixes the chart series formulas on ZONE_CALCULATORS so that the category/value
#     ranges keep pointing at the right sheet once it is renamed.
#  2. Renames the four worksheets so the underscores become spaces:
#       ZONE_CALCULATORS        -> ZONE CALCULATORS
#       RESOURCE_MGR            -> RESOURCE MGR
#       UPLOAD_READY_PRODUCTION -> UPLOAD READY PRODUCTION
#       CROSS_REFERENCE         -> CROSS REFERENCE
#  3. Updates the "Total Payroll Forecast" figure on CROSS REFERENCE (cell B22) from 0
#     to 80000 (CPO/payroll inflation fix).

$wb = $excel.ActiveWorkbook

$oldZone  = "ZONE_CALCULATORS"
$newZone  = "ZONE CALCULATORS"
$oldRes   = "RESOURCE_MGR"
$newRes   = "RESOURCE MGR"
$oldUp    = "UPLOAD_READY_PRODUCTION"
$newUp    = "UPLOAD READY PRODUCTION"
$oldCross = "CROSS_REFERENCE"
$newCross = "CROSS REFERENCE"

$wsZone  = $wb.Worksheets.Item($oldZone)
$wsRes   = $wb.Worksheets.Item($oldRes)
$wsUp    = $wb.Worksheets.Item($oldUp)
$wsCross = $wb.Worksheets.Item($oldCross)

# --- Step 1: fix the chart series formulas that live on ZONE_CALCULATORS so the chart
#             keeps plotting the right ranges after the sheet gets its new name ---
$chartObjs = $wsZone.ChartObjects()
for ($i = 1; $i -le $chartObjs.Count; $i++) {
    $co = $chartObjs.Item($i)
    $chart = $co.Chart
    $sc = $chart.SeriesCollection()
    for ($j = 1; $j -le $sc.Count; $j++) {
        $ser = $sc.Item($j)
        $f = $ser.Formula
        if ($f -like "*$oldZone*") {
            $ser.Formula = $f.Replace($oldZone, "'$newZone'")
        }
    }
}

# --- Step 2: rename the sheets (spaces instead of underscores) ---
$wsZone.Name = $newZone
$wsRes.Name = $newRes
$wsUp.Name = $newUp
$wsCross.Name = $newCross

# --- Step 3: bump the Total Payroll Forecast figure on CROSS REFERENCE ---
$wsCross.Range("B22").Value = 80000
